$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 9.226618575922256
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 19.51936550083139

# Row 3
$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 1766.335244827366
$ws.Range("D3").Value = 2938.103010863317
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 4712.425297809124

# Row 4
$ws.Range("B4").Value = 0.06328177979961902
$ws.Range("C4").Value = 0.05231270169004087
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 0.768386970581898

# Row 5
$ws.Range("B5").Value = 0.3464964993005633
$ws.Range("C5").Value = 0.3375848360084654
$ws.Range("D5").Value = 0.1529057820181812
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("G5").Value = 1.336873824401267
